$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly-refreshed "Fruta / hortaliza" price feed: every
# refresh appends newly observed weekly records into the existing date-sorted
# table and re-sorts by date, which shows up as a mix of value changes plus
# 3 brand-new rows. Net effect here: the data block grows from 34 to 37 rows
# (dimension A1:T34 -> A1:T37), with rows 12-34 getting updated contents and
# 3 new rows appended.

# Insert 3 new blank rows at the bottom of the data block (old rows 35-37)
# so the sheet grows to the new size; we then rewrite rows 12-37 in full
# below with the refreshed, re-sorted dataset.
$ws.Rows("35:37").Insert()

# Values shared by every data row in this block (constant across the table)
$colA = 8
$colB = "Terminal La Palmera de La Serena"
$colC = "Coquimbo"
$colE = 4
$colF = "Fruta"
$colG = 100108
$colH = "Tropicales y subtropicales"
$colI = 100108007
$colJ = "Coco"
$colK = "Sin especificar"
$colL = "Primera"
$colQ = '$/malla 20 unidades'
$colR = "Perú"
$colT = 20

# Per-row data: @(row, Fecha, Volumen, PrecioMin, PrecioMax, PrecioPromPonderado, PrecioPorKg)
$rowsData = @(
    @(12, 44810, 100, 27000, 28000, 27500, 1375),
    @(13, 44350, 160, 19000, 20000, 19500, 975),
    @(14, 44407, 160, 20000, 21000, 20500, 1025),
    @(15, 44784, 160, 27000, 28000, 27500, 1375),
    @(16, 44776, 160, 23000, 24000, 23500, 1175),
    @(17, 44466, 100, 20000, 21000, 20500, 1025),
    @(18, 44427, 200, 20000, 21000, 20500, 1025),
    @(19, 44410, 200, 20000, 21000, 20500, 1025),
    @(20, 44326, 160, 19500, 20000, 19750, 988),
    @(21, 44418, 200, 20000, 21000, 20500, 1025),
    @(22, 44441, 160, 20000, 21000, 20500, 1025),
    @(23, 44335, 200, 19000, 20000, 19500, 975),
    @(24, 44448, 100, 20000, 21000, 20500, 1025),
    @(25, 44809, 60, 27000, 28000, 27500, 1375),
    @(26, 44315, 100, 20000, 21000, 20500, 1025),
    @(27, 44778, 100, 23000, 24000, 23500, 1175),
    @(28, 44431, 160, 21000, 22000, 21500, 1075),
    @(29, 44434, 100, 20000, 21000, 20500, 1025),
    @(30, 44462, 100, 19500, 20000, 19750, 988),
    @(31, 44442, 140, 20000, 21000, 20500, 1025),
    @(32, 44781, 160, 23000, 24000, 23500, 1175),
    @(33, 44336, 100, 19500, 20000, 19750, 988),
    @(34, 44343, 100, 19500, 20000, 19750, 988),
    @(35, 44365, 100, 20000, 21000, 20500, 1025),
    @(36, 44782, 200, 23500, 24000, 23750, 1188),
    @(37, 44435, 260, 20000, 22000, 21115, 1056)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ
    $ws.Cells.Item($r, 11).Value = $colK
    $ws.Cells.Item($r, 12).Value = $colL
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 17).Value = $colQ
    $ws.Cells.Item($r, 18).Value = $colR
    $ws.Cells.Item($r, 19).Value = $row[6]
    $ws.Cells.Item($r, 20).Value = $colT
}
